$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$cl21 = $m.CustomLayouts.Item(21)
$cl21.Name = "video meldingXXX"
